$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.128.95'
$ws.Range('E2').Value = '  -2.32%  '
$ws.Range('D3').Value = '1.564.74'
$ws.Range('E3').Value = '  -1.92%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = '''206.38'
$ws.Range('E5').Value = '  -1.40%  '
$ws.Range('E6').Value = '  -3.23%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').Value = '''22.25'
$ws.Range('E8').Value = '  -0.52%  '
$ws.Range('D9').Value = '''0.247'
$ws.Range('E9').Value = '  -2.58%  '
$ws.Range('D11').Value = '''0.0862'
$ws.Range('E11').Value = '  -0.69%  '
$ws.Range('D12').Value = '1.785.83'
$ws.Range('E12').Value = '  -2.04%  '
$ws.Range('D13').Value = '1.571.28'
$ws.Range('E13').Value = '  -2.02%  '
$ws.Range('D14').Value = '''3.77'
$ws.Range('E14').Value = '  -2.28%  '
$ws.Range('E15').Value = '  -3.23%  '
$ws.Range('D16').Value = '''62.99'
$ws.Range('E16').Value = '  -0.79%  '
$ws.Range('D17').Value = '27.127.39'
$ws.Range('E17').Value = '  -2.34%  '
$ws.Range('D18').Value = '''214.88'
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range('E19').Value = '  -1.93%  '
$ws.Range('E20').Value = '  -1.94%  '
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('E23').Value = '  -4.35%  '
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('D25').Value = '''151.78'
$ws.Range('E25').Value = '  -1.47%  '
$ws.Range('E26').Value = '  -8.01%  '
$ws.Range('E27').Value = '  -1.62%  '
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('E29').Value = '  -1.65%  '
$ws.Range('E30').Value = '  -2.13%  '
$ws.Range('D31').Value = '''0.0462'
$ws.Range('E31').Value = '  -2.54%  '
$ws.Range('D32').Value = '''3.16'
$ws.Range('E32').Value = '  -2.21%  '
$ws.Range('D33').Value = '1.393.05'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('E36').Value = '  -2.11%  '
$ws.Range('E37').Value = '  -3.64%  '
$ws.Range('E38').Value = '  -2.12%  '
$ws.Range('E39').Value = '  -2.12%  '
$ws.Range('D40').Value = '''0.516'
$ws.Range('E40').Value = '  -3.89%  '
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('D42').Value = '''0.991'
$ws.Range('E42').Value = '  +1.97%  '
$ws.Range('E43').Value = '  +2.28%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''5.29'
$ws.Range('E44').Value = '  +0.68%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '''63.49'
$ws.Range('E45').Value = '  -1.70%  '
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('D47').Value = '1.699.11'
$ws.Range('E47').Value = '  -2.03%  '
$ws.Range('D48').Value = '''85.52'
$ws.Range('E48').Value = '  -1.28%  '
$ws.Range('D49').Value = '0.0₇0988'
$ws.Range('E49').Value = '  -2.51%  '
$ws.Range('E50').Value = '  -0.93%  '
$ws.Range('D51').Value = '''0.0946'
$ws.Range('E51').Value = '  -2.17%  '
